# Fixed naive component forecaster bug - Presentation state 11.02.
#
# For each data row (2..20), a new forecast-error value (for the newly
# fixed horizon Q0) is inserted into column B, and all the existing
# values that were in B..K are shifted one column to the right
# (B->C, C->D, ..., J->K). If the row was already "full" (had a value
# in K), that trailing value is dropped off the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to insert into column B for each row (2..20).
$newValues = @{
    2  = 0.2433086034063205
    3  = -0.115952868393072
    4  = 0.1459581181302581
    5  = -0.08819670345554087
    6  = 0.3847923593882046
    7  = 0.03766489642184559
    8  = 0.1769978556124878
    9  = 0.1302808926112106
    10 = -0.1944981035472806
    11 = -0.3817857436446591
    12 = -0.0000002375649628613696
    13 = -0.0555296279974082
    14 = 0.0000003829984367986761
    15 = -0.0000001604754923945073
    16 = 0.009398958989038461
    17 = -0.07651818316594991
    18 = 0.0000002770877186031306
    19 = 0.2010531357750048
    20 = -0.2003621554241067
}

$firstCol = 2   # column B
$lastCol  = 11  # column K

for ($r = 2; $r -le 20; $r++) {

    # Read the existing row values (columns B..K) before modifying anything.
    $existing = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $existing += , ($ws.Cells.Item($r, $c).Value())
    }

    # Build the shifted row: new value first, then the previous values,
    # truncated so the row never grows past column K.
    $shifted = @(, $newValues[$r])
    foreach ($val in $existing) {
        $shifted += , $val
    }

    # Write the new row values back out, column by column.
    for ($i = 0; $i -lt $existing.Count; $i++) {
        $col = $firstCol + $i
        $newVal = $shifted[$i]
        if ($null -eq $newVal) {
            $ws.Cells.Item($r, $col).Value = ""
        } else {
            $ws.Cells.Item($r, $col).Value = $newVal
        }
    }
}
